# Workbook "default-ts-ml-Tx_var.xlsx" - convert transmitter incidence angle
# (Tx_th, measured from zenith) into transmitter elevation angle (Tx_el,
# measured from the horizon) on the "Dynamic" sheet.
#   Tx_el = 90 - Tx_th
# This also renames the column B header from "Tx_th (deg)" to "Tx_el (deg)".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dynamic")

# Rename header (column B, row 1): "Tx_th (deg)" -> "Tx_el (deg)"
$ws.Range("B1").Value2 = "Tx_el (deg)"

# Convert every data value in column B (rows 2-406) from incidence angle to
# elevation angle: new = 90 - old
$lastRow = 406
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $current = $cell.Value2
    if ($current -ne $null) {
        $cell.Value2 = 90 - $current
    }
}

# Restore the last-used selection recorded in the sheet view
[void]$ws.Range("P13").Select()
